# Generate Report for Handback
# The "b449e6a1-ccaf-46f7-9753-88a10d2650d6.md" file has been handed back
# (in sync with en-US) for both zh-cn and de-de locales, so update the
# status/handback-datetime cells across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $status
$zhcn.Range("G3").Value = "2016-03-04 03:28:46"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $status
$dede.Range("G3").Value = "2016-03-04 03:29:11"
